$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.332.16'
$ws.Range('E2').Value = '  +2.26%  '
$ws.Range('D3').Value = '3.150.02'
$ws.Range('E3').Value = '  +2.77%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '537.34'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +3.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.75'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.41%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.515'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +10.50%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.33'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.57%  '
$ws.Range('E10').Value = '  +3.39%  '
$ws.Range('E11').Value = '  +5.49%  '
$ws.Range('E12').Value = '  +2.80%  '
$ws.Range('D13').Value = '3.692.03'
$ws.Range('E13').Value = '  +2.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.03'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.80%  '
$ws.Range('E15').Value = '  +6.00%  '
$ws.Range('D16').Value = '58.377.53'
$ws.Range('E16').Value = '  +2.20%  '
$ws.Range('D17').Value = '3.152.06'
$ws.Range('E17').Value = '  +2.82%  '
$ws.Range('E18').Value = '  +6.42%  '
$ws.Range('E19').Value = '  +4.85%  '
$ws.Range('E20').Value = '  +5.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '377.18'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +8.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.74'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.50%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.51'
$ws.Range('D24').ClearFormats()
$ws.Range('E25').Value = '  +4.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.168'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.21%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.00'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +12.20%  '
$ws.Range('E29').Value = '  +3.19%  '
$ws.Range('E30').Value = '  +2.90%  '
$ws.Range('E31').Value = '  +7.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.83'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +4.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.17'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +7.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.18'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +5.24%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '161.49'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.76%  '
$ws.Range('E36').Value = '  +4.90%  '
$ws.Range('E37').Value = '  +10.25%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.57'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.64%  '
$ws.Range('E39').Value = '  +8.10%  '
$ws.Range('D40').Value = '2.648.94'
$ws.Range('E40').Value = '  +10.27%  '
$ws.Range('E41').Value = '  +4.23%  '
$ws.Range('E42').Value = '  +5.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '38.65'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +6.38%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.704'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.90%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0276'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +6.20%  '
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('E47').Value = '  +12.54%  '
$ws.Range('E48').Value = '  +4.69%  '
$ws.Range('E49').Value = '  +5.35%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.25'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +4.42%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.752'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.65%  '
